$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: row 3 (bfbf0c67... file) gets a refreshed handback report
# (new handoff/handback datetimes recorded for this handback run).
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-17 06:32:00"
$wsZh.Range("H3").Value = "2016-03-17 06:32:41"

# "de-de" sheet: row 3 (bfbf0c67... file) gets the same treatment.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-17 06:32:08"
$wsDe.Range("H3").Value = "2016-03-17 06:32:54"
